$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '42.882.37'
$ws.Range("E2").Value = '  -5.30%  '
$ws.Range("D3").Value = '2.208.74'
$ws.Range("E3").Value = '  -6.70%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '314.89'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.68%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '97.52'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -9.87%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.581'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -7.59%  '
$ws.Range("E8").Value = '  +0.07%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.556'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -9.76%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '36.42'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -11.23%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '54.25'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.97%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.0824'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -10.36%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '7.71'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -8.84%  '
$ws.Range("E14").Value = '  -3.99%  '
$ws.Range("E15").Value = '  -12.07%  '
$ws.Range("D16").Value = '2.546.97'
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '14.05'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -7.61%  '
$ws.Range("D18").Value = '2.208.48'
$ws.Range("E18").Value = '  -6.63%  '
$ws.Range("D19").Value = '42.741.86'
$ws.Range("E19").Value = '  -5.54%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '14.54'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("E21").Value = '  -9.95%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.36'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -12.79%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '65.15'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -10.93%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '3.16'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -8.94%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '235.70'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -9.51%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.11'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -8.59%  '
$ws.Range("E27").Value = '  +0.24%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '10.00'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -10.26%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '6.19'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -15.44%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '20.38'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -8.78%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.0875'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -9.69%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '33.63'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -10.53%  '
$ws.Range("E34").Value = '  -8.56%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.77'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -6.08%  '
$ws.Range("E36").Value = '  +7.52%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.98'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +13.68%  '
$ws.Range("E38").Value = '  -6.78%  '
$ws.Range("E39").Value = '  -7.97%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.102'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -12.82%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '3.68'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -5.97%  '
$ws.Range("E42").Value = '  -8.71%  '
$ws.Range("D43").Value = '1.860.59'
$ws.Range("E43").Value = '  +11.57%  '
$ws.Range("E44").Value = '  +0.03%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '12.21'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -5.48%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '88.56'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -11.42%  '
$ws.Range("E47").Value = '  -11.03%  '
$ws.Range("E48").Value = '  -2.42%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '75.58'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -6.35%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '59.75'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -14.01%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '8.61'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -6.31%  '
